$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell E1
$ws.Range("E1").Value = "Time(ms)"

# Update data rows (rows 2-26 correspond to generations 0-24)
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = "pred wins"
$ws.Cells.Item(2, 3).Value = 0
$ws.Cells.Item(2, 4).Value = 8
$ws.Cells.Item(2, 5).Value = 8273.701429367065

$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = "prey wins"
$ws.Cells.Item(3, 3).Value = 38
$ws.Cells.Item(3, 4).Value = 10
$ws.Cells.Item(3, 5).Value = 2480.821132659912

$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = "prey wins"
$ws.Cells.Item(4, 3).Value = 37
$ws.Cells.Item(4, 4).Value = 11
$ws.Cells.Item(4, 5).Value = 2461.47084236145

$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(5, 2).Value = "prey wins"
$ws.Cells.Item(5, 3).Value = 0
$ws.Cells.Item(5, 4).Value = 15
$ws.Cells.Item(5, 5).Value = 2444.170713424683

$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 2).Value = "prey wins"
$ws.Cells.Item(6, 3).Value = 22
$ws.Cells.Item(6, 4).Value = 13
$ws.Cells.Item(6, 5).Value = 2469.356060028076

$ws.Cells.Item(7, 1).Value = 5
$ws.Cells.Item(7, 2).Value = "prey wins"
$ws.Cells.Item(7, 3).Value = 50
$ws.Cells.Item(7, 4).Value = 12
$ws.Cells.Item(7, 5).Value = 2513.115882873535

$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = "prey wins"
$ws.Cells.Item(8, 3).Value = 44
$ws.Cells.Item(8, 4).Value = 14
$ws.Cells.Item(8, 5).Value = 2753.872632980347

$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "prey wins"
$ws.Cells.Item(9, 3).Value = 49
$ws.Cells.Item(9, 4).Value = 18
$ws.Cells.Item(9, 5).Value = 2481.059789657593

$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = "prey wins"
$ws.Cells.Item(10, 3).Value = 59
$ws.Cells.Item(10, 4).Value = 16
$ws.Cells.Item(10, 5).Value = 2607.418537139893

$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = "prey wins"
$ws.Cells.Item(11, 3).Value = 56
$ws.Cells.Item(11, 4).Value = 17
$ws.Cells.Item(11, 5).Value = 2484.508514404297

$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 2).Value = "prey wins"
$ws.Cells.Item(12, 3).Value = 56
$ws.Cells.Item(12, 4).Value = 18
$ws.Cells.Item(12, 5).Value = 2673.427581787109

$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 2).Value = "prey wins"
$ws.Cells.Item(13, 3).Value = 62
$ws.Cells.Item(13, 4).Value = 19
$ws.Cells.Item(13, 5).Value = 2567.131757736206

$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(14, 2).Value = "prey wins"
$ws.Cells.Item(14, 3).Value = 63
$ws.Cells.Item(14, 4).Value = 19
$ws.Cells.Item(14, 5).Value = 2536.525726318359

$ws.Cells.Item(15, 1).Value = 13
$ws.Cells.Item(15, 2).Value = "prey wins"
$ws.Cells.Item(15, 3).Value = 61
$ws.Cells.Item(15, 4).Value = 19
$ws.Cells.Item(15, 5).Value = 2523.674011230469

$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "prey wins"
$ws.Cells.Item(16, 3).Value = 35
$ws.Cells.Item(16, 4).Value = 6
$ws.Cells.Item(16, 5).Value = 4960.159301757812

$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = "prey wins"
$ws.Cells.Item(17, 3).Value = 47
$ws.Cells.Item(17, 4).Value = 12
$ws.Cells.Item(17, 5).Value = 2468.917608261108

$ws.Cells.Item(18, 1).Value = 16
$ws.Cells.Item(18, 2).Value = "prey wins"
$ws.Cells.Item(18, 3).Value = 61
$ws.Cells.Item(18, 4).Value = 12
$ws.Cells.Item(18, 5).Value = 2480.092763900757

$ws.Cells.Item(19, 1).Value = 17
$ws.Cells.Item(19, 2).Value = "prey wins"
$ws.Cells.Item(19, 3).Value = 38
$ws.Cells.Item(19, 4).Value = 14
$ws.Cells.Item(19, 5).Value = 2477.5230884552

$ws.Cells.Item(20, 1).Value = 18
$ws.Cells.Item(20, 2).Value = "prey wins"
$ws.Cells.Item(20, 3).Value = 56
$ws.Cells.Item(20, 4).Value = 15
$ws.Cells.Item(20, 5).Value = 2709.322214126587

$ws.Cells.Item(21, 1).Value = 19
$ws.Cells.Item(21, 2).Value = "prey wins"
$ws.Cells.Item(21, 3).Value = 48
$ws.Cells.Item(21, 4).Value = 16
$ws.Cells.Item(21, 5).Value = 2655.769348144531

$ws.Cells.Item(22, 1).Value = 20
$ws.Cells.Item(22, 2).Value = "prey wins"
$ws.Cells.Item(22, 3).Value = 43
$ws.Cells.Item(22, 4).Value = 16
$ws.Cells.Item(22, 5).Value = 2514.684915542603

$ws.Cells.Item(23, 1).Value = 21
$ws.Cells.Item(23, 2).Value = "prey wins"
$ws.Cells.Item(23, 3).Value = 49
$ws.Cells.Item(23, 4).Value = 18
$ws.Cells.Item(23, 5).Value = 2538.419485092163

$ws.Cells.Item(24, 1).Value = 22
$ws.Cells.Item(24, 2).Value = "prey wins"
$ws.Cells.Item(24, 3).Value = 52
$ws.Cells.Item(24, 4).Value = 17
$ws.Cells.Item(24, 5).Value = 2509.968280792236

$ws.Cells.Item(25, 1).Value = 23
$ws.Cells.Item(25, 2).Value = "prey wins"
$ws.Cells.Item(25, 3).Value = 42
$ws.Cells.Item(25, 4).Value = 18
$ws.Cells.Item(25, 5).Value = 2498.693466186523

$ws.Cells.Item(26, 1).Value = 24
$ws.Cells.Item(26, 2).Value = "prey wins"
$ws.Cells.Item(26, 3).Value = 61
$ws.Cells.Item(26, 4).Value = 18
$ws.Cells.Item(26, 5).Value = 2559.166431427002
